$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (date 2021-11-10 / serial 44510) was inserted
# as row 4, pushing the previously existing rows 4-7 down to rows 5-8.
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44510
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 300000000
$ws.Range("G4").Value = "Espárragos"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 950
$ws.Range("N4").Value = "`$/kilo"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 950
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
